$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix time-calculation error: total hours value in B1 was wrong (60 -> 36)
$ws.Range("B1").Value = 36

# Move the active selection to A16 (matches post-edit cursor position)
$null = $ws.Range("A16").Select()
